$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Farmers Database" - add a new "Coconut" crop column (H) and update
# a couple of farmer records.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Farmers Database")

# New header for the Coconut crop column (copy the header formatting from
# the neighbouring "Cashewnuts" header cell, then set the text).
$ws1.Range("G1").Copy()
$ws1.Range("H1").PasteSpecial(-4122)
$ws1.Range("H1").Value = "Coconut"

# Arvi (row 2) logged 5 units of Banana instead of the earlier 1000 of Rice.
$ws1.Range("E2").Value = 5
$ws1.Range("F2").Value = 0
$ws1.Range("H2").Value = 0

# New Coconut column defaults to 0 for the other existing farmers.
$ws1.Range("H3").Value = 0
$ws1.Range("H4").Value = 0

# Arvind's (row 5) details were corrected/simplified.
$ws1.Range("B5").Value = 10
$ws1.Range("C5").Value = "A"
$ws1.Range("D5").Value = 10
$ws1.Range("H5").Value = 0

# ---------------------------------------------------------------------------
# Sheet "Farmers Log" - the old submissions were cleared out and replaced by
# two fresh log entries.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Farmers Log")

# Drop the third (oldest-retained) log row entirely; the remaining two rows
# get overwritten below with the new submissions.
$ws2.Rows.Item(4).Delete()

$ws2.Range("A2").Value = "2017-10-02 10:38:48"
$ws2.Range("B2").Value = "Shan"
$ws2.Range("C2").Value = 12343
$ws2.Range("D2").Value = "Rice"
$ws2.Range("E2").Value = 500

$ws2.Range("A3").Value = "2017-10-02 11:31:09"
$ws2.Range("B3").Value = "Arvi"
$ws2.Range("C3").Value = 12341
$ws2.Range("D3").Value = "Banana"
$ws2.Range("E3").Value = 10

# ---------------------------------------------------------------------------
# Sheet "Companies Log" - the Nestle purchase record was updated with a new
# timestamp and revised quantity/price.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Companies Log")

$ws3.Range("A2").Value = "2017-10-02 11:35:05"
$ws3.Range("D2").Value = 5
$ws3.Range("E2").Value = 20500
